$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Scratch cell outside the used range (A1:G51) that carries the sheet's
# default (unstyled) format; used below to strip the temporary text-number
# format back off each edited cell after we write its new value.
$blank = $ws.Range("H1")

function Set-TextValue($addr, $val) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $blank.Copy() | Out-Null
    $cell.PasteSpecial(-4122) | Out-Null  # xlPasteFormats: restore original (General) style
}

Set-TextValue "D2" "278.17"
Set-TextValue "E2" "0.72%"
Set-TextValue "D3" "27.23"
Set-TextValue "E3" "2.38%"
Set-TextValue "D4" "4.865"
Set-TextValue "E4" "-0.17%"
Set-TextValue "D5" "0.06421"
Set-TextValue "E5" "1.50%"
Set-TextValue "D6" "7.020"
Set-TextValue "E6" "1.47%"
Set-TextValue "D7" "1.205"
Set-TextValue "E7" "-7.55%"
Set-TextValue "D8" "0.8863"
Set-TextValue "E8" "1.28%"
Set-TextValue "D9" "0.1552"
Set-TextValue "E9" "-0.26%"
Set-TextValue "D10" "0.05106"
Set-TextValue "E10" "1.97%"
Set-TextValue "D11" "0.07520"
Set-TextValue "E11" "0.65%"
Set-TextValue "D12" "0.02886"
Set-TextValue "E12" "-3.07%"
Set-TextValue "D13" "0.08964"
Set-TextValue "E13" "-1.04%"
Set-TextValue "D14" "0.001565"
Set-TextValue "E14" "-0.50%"
Set-TextValue "D15" "0.0006393"
Set-TextValue "E15" "1.36%"
Set-TextValue "D16" "0.006085"
Set-TextValue "E16" "0.48%"
Set-TextValue "D17" "3.478"
Set-TextValue "E17" "0.92%"
Set-TextValue "D18" "3.306"
Set-TextValue "E18" "-0.41%"
Set-TextValue "E21" "-0.13%"
Set-TextValue "D22" "3.905"
Set-TextValue "E22" "0.19%"
Set-TextValue "D23" "0.04417"
Set-TextValue "E23" "1.10%"
Set-TextValue "E25" "0.45%"
Set-TextValue "D26" "0.003880"
Set-TextValue "E26" "-7.86%"
Set-TextValue "E28" "-1.62%"
Set-TextValue "E29" "1.73%"
Set-TextValue "D40" "0.04127"
Set-TextValue "E40" "0.63%"
Set-TextValue "D41" "0.006790"
Set-TextValue "E41" "-2.75%"
Set-TextValue "E42" "0.08%"
Set-TextValue "D43" "0.001921"
Set-TextValue "E43" "-9.39%"
Set-TextValue "E44" "3.12%"
Set-TextValue "D45" "0.00005315"
Set-TextValue "E45" "0.24%"
Set-TextValue "E46" "13.10%"
Set-TextValue "E47" "-7.42%"
